$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Fitness column (C2:C12) to the new constant value 4171
$ws.Range("C2:C12").Value = 4171
